# Update NATMI LR-pair sheet (Grn-Tnfrsf1a) with recomputed TPM-derived values.
# Ligand (G/H -> I/J) and Receptor (M/N -> O/P) expression values change per
# sending/target cluster; Q/R/S/T (edge weights/specificities) are the derived
# products updated to match.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 43.34730933333333
$ws.Range("H2").Value = 130.041928
$ws.Range("I2").Value = 0.04273139820300816
$ws.Range("J2").Value = 0.04273139820300816
$ws.Range("M2").Value = 51.402972
$ws.Range("N2").Value = 154.208916
$ws.Range("O2").Value = 0.2478201393026385
$ws.Range("P2").Value = 0.2478201393026385
$ws.Range("Q2").Value = 2228.180527936672
$ws.Range("R2").Value = 20053.62475143004
$ws.Range("S2").Value = 0.010589701055266
$ws.Range("T2").Value = 0.010589701055266

# Row 3
$ws.Range("G3").Value = 43.34730933333333
$ws.Range("H3").Value = 130.041928
$ws.Range("I3").Value = 0.04273139820300816
$ws.Range("J3").Value = 0.04273139820300816
$ws.Range("O3").Value = 0.2505855219821376
$ws.Range("P3").Value = 0.2505855219821376
$ws.Range("Q3").Value = 2253.044414528344
$ws.Range("R3").Value = 20277.3997307551
$ws.Range("S3").Value = 0.01070786972372738
$ws.Range("T3").Value = 0.01070786972372737

# Row 4
$ws.Range("G4").Value = 43.34730933333333
$ws.Range("H4").Value = 130.041928
$ws.Range("I4").Value = 0.04273139820300816
$ws.Range("J4").Value = 0.04273139820300816
$ws.Range("M4").Value = 59.18509700000001
$ws.Range("N4").Value = 177.555291
$ws.Range("O4").Value = 0.2853387345614992
$ws.Range("P4").Value = 0.2853387345614991
$ws.Range("Q4").Value = 2565.514707582339
$ws.Range("R4").Value = 23089.63236824105
$ws.Range("S4").Value = 0.01219292308928987
$ws.Range("T4").Value = 0.01219292308928986

# Row 5
$ws.Range("G5").Value = 43.34730933333333
$ws.Range("H5").Value = 130.041928
$ws.Range("I5").Value = 0.04273139820300816
$ws.Range("J5").Value = 0.04273139820300816
$ws.Range("M5").Value = 8.716382666666666
$ws.Range("N5").Value = 26.149148
$ws.Range("O5").Value = 0.04202276799614693
$ws.Range("P5").Value = 0.04202276799614692
$ws.Range("Q5").Value = 377.8317357197048
$ws.Range("R5").Value = 3400.485621477343
$ws.Range("S5").Value = 0.001795691632835982
$ws.Range("T5").Value = 0.001795691632835981

# Row 6
$ws.Range("G6").Value = 43.34730933333333
$ws.Range("H6").Value = 130.041928
$ws.Range("I6").Value = 0.04273139820300816
$ws.Range("J6").Value = 0.04273139820300816
$ws.Range("M6").Value = 36.13945833333333
$ws.Range("N6").Value = 108.418375
$ws.Range("O6").Value = 0.1742328361575779
$ws.Range("P6").Value = 0.1742328361575779
$ws.Range("Q6").Value = 1566.548279514111
$ws.Range("R6").Value = 14098.934515627
$ws.Range("S6").Value = 0.007445212701888941
$ws.Range("T6").Value = 0.00744521270188894

# Row 7
$ws.Range("I7").Value = 0.0889365509391893
$ws.Range("J7").Value = 0.08893655093918929
$ws.Range("M7").Value = 51.402972
$ws.Range("N7").Value = 154.208916
$ws.Range("O7").Value = 0.2478201393026385
$ws.Range("P7").Value = 0.2478201393026385
$ws.Range("Q7").Value = 4637.496065143948
$ws.Range("R7").Value = 41737.46458629553
$ws.Range("S7").Value = 0.02204026844284609
$ws.Range("T7").Value = 0.02204026844284609

# Row 8
$ws.Range("I8").Value = 0.0889365509391893
$ws.Range("J8").Value = 0.08893655093918929
$ws.Range("O8").Value = 0.2505855219821376
$ws.Range("P8").Value = 0.2505855219821376
$ws.Range("S8").Value = 0.02228621204038772
$ws.Range("T8").Value = 0.02228621204038771

# Row 9
$ws.Range("I9").Value = 0.0889365509391893
$ws.Range("J9").Value = 0.08893655093918929
$ws.Range("M9").Value = 59.18509700000001
$ws.Range("N9").Value = 177.555291
$ws.Range("O9").Value = 0.2853387345614992
$ws.Range("P9").Value = 0.2853387345614991
$ws.Range("Q9").Value = 5339.58726068724
$ws.Range("R9").Value = 48056.28534618516
$ws.Range("S9").Value = 0.02537704290125259
$ws.Range("T9").Value = 0.02537704290125258

# Row 10
$ws.Range("I10").Value = 0.0889365509391893
$ws.Range("J10").Value = 0.08893655093918929
$ws.Range("M10").Value = 8.716382666666666
$ws.Range("N10").Value = 26.149148
$ws.Range("O10").Value = 0.04202276799614693
$ws.Range("P10").Value = 0.04202276799614692
$ws.Range("Q10").Value = 786.3784669679328
$ws.Range("R10").Value = 7077.406202711395
$ws.Range("S10").Value = 0.003737360046495056
$ws.Range("T10").Value = 0.003737360046495054

# Row 11
$ws.Range("I11").Value = 0.0889365509391893
$ws.Range("J11").Value = 0.08893655093918929
$ws.Range("M11").Value = 36.13945833333333
$ws.Range("N11").Value = 108.418375
$ws.Range("O11").Value = 0.1742328361575779
$ws.Range("P11").Value = 0.1742328361575779
$ws.Range("Q11").Value = 3260.445637603736
$ws.Range("R11").Value = 29344.01073843363
$ws.Range("S11").Value = 0.01549566750820785
$ws.Range("T11").Value = 0.01549566750820785

# Row 12
$ws.Range("G12").Value = 394.701121
$ws.Range("H12").Value = 1184.103363
$ws.Range("I12").Value = 0.3890929110023202
$ws.Range("J12").Value = 0.3890929110023201
$ws.Range("M12").Value = 51.402972
$ws.Range("N12").Value = 154.208916
$ws.Range("O12").Value = 0.2478201393026385
$ws.Range("P12").Value = 0.2478201393026385
$ws.Range("Q12").Value = 20288.81067113161
$ws.Range("R12").Value = 182599.2960401845
$ws.Range("S12").Value = 0.0964250594062641
$ws.Range("T12").Value = 0.09642505940626407

# Row 13
$ws.Range("G13").Value = 394.701121
$ws.Range("H13").Value = 1184.103363
$ws.Range("I13").Value = 0.3890929110023202
$ws.Range("J13").Value = 0.3890929110023201
$ws.Range("O13").Value = 0.2505855219821376
$ws.Range("P13").Value = 0.2505855219821376
$ws.Range("Q13").Value = 20515.21005003385
$ws.Range("R13").Value = 184636.8904503046
$ws.Range("S13").Value = 0.09750105020306581
$ws.Range("T13").Value = 0.09750105020306578

# Row 14
$ws.Range("G14").Value = 394.701121
$ws.Range("H14").Value = 1184.103363
$ws.Range("I14").Value = 0.3890929110023202
$ws.Range("J14").Value = 0.3890929110023201
$ws.Range("M14").Value = 59.18509700000001
$ws.Range("N14").Value = 177.555291
$ws.Range("O14").Value = 0.2853387345614992
$ws.Range("P14").Value = 0.2853387345614991
$ws.Range("Q14").Value = 23360.42413239374
$ws.Range("R14").Value = 210243.8171915436
$ws.Range("S14").Value = 0.1110232788522521
$ws.Range("T14").Value = 0.111023278852252

# Row 15
$ws.Range("G15").Value = 394.701121
$ws.Range("H15").Value = 1184.103363
$ws.Range("I15").Value = 0.3890929110023202
$ws.Range("J15").Value = 0.3890929110023201
$ws.Range("M15").Value = 8.716382666666666
$ws.Range("N15").Value = 26.149148
$ws.Range("O15").Value = 0.04202276799614693
$ws.Range("P15").Value = 0.04202276799614692
$ws.Range("Q15").Value = 3440.366009598302
$ws.Range("R15").Value = 30963.29408638472
$ws.Range("S15").Value = 0.01635076112799595
$ws.Range("T15").Value = 0.01635076112799594

# Row 16
$ws.Range("G16").Value = 394.701121
$ws.Range("H16").Value = 1184.103363
$ws.Range("I16").Value = 0.3890929110023202
$ws.Range("J16").Value = 0.3890929110023201
$ws.Range("M16").Value = 36.13945833333333
$ws.Range("N16").Value = 108.418375
$ws.Range("O16").Value = 0.1742328361575779
$ws.Range("P16").Value = 0.1742328361575779
$ws.Range("Q16").Value = 14264.28471649946
$ws.Range("R16").Value = 128378.5624484951
$ws.Range("S16").Value = 0.06779276141274231
$ws.Range("T16").Value = 0.06779276141274228

# Row 17
$ws.Range("G17").Value = 7.804371333333333
$ws.Range("H17").Value = 23.413114
$ws.Range("I17").Value = 0.007693480963358413
$ws.Range("J17").Value = 0.007693480963358412
$ws.Range("M17").Value = 51.402972
$ws.Range("N17").Value = 154.208916
$ws.Range("O17").Value = 0.2478201393026385
$ws.Range("P17").Value = 0.2478201393026385
$ws.Range("Q17").Value = 401.167881124936
$ws.Range("R17").Value = 3610.510930124424
$ws.Range("S17").Value = 0.001906599524061679
$ws.Range("T17").Value = 0.001906599524061679

# Row 18
$ws.Range("G18").Value = 7.804371333333333
$ws.Range("H18").Value = 23.413114
$ws.Range("I18").Value = 0.007693480963358413
$ws.Range("J18").Value = 0.007693480963358412
$ws.Range("O18").Value = 0.2505855219821376
$ws.Range("P18").Value = 0.2505855219821376
$ws.Range("Q18").Value = 405.6444451086221
$ws.Range("R18").Value = 3650.800005977598
$ws.Range("S18").Value = 0.001927874943062806
$ws.Range("T18").Value = 0.001927874943062806

# Row 19
$ws.Range("G19").Value = 7.804371333333333
$ws.Range("H19").Value = 23.413114
$ws.Range("I19").Value = 0.007693480963358413
$ws.Range("J19").Value = 0.007693480963358412
$ws.Range("M19").Value = 59.18509700000001
$ws.Range("N19").Value = 177.555291
$ws.Range("O19").Value = 0.2853387345614992
$ws.Range("P19").Value = 0.2853387345614991
$ws.Range("Q19").Value = 461.9024743873527
$ws.Range("R19").Value = 4157.122269486174
$ws.Range("S19").Value = 0.002195248122457673
$ws.Range("T19").Value = 0.002195248122457672

# Row 20
$ws.Range("G20").Value = 7.804371333333333
$ws.Range("H20").Value = 23.413114
$ws.Range("I20").Value = 0.007693480963358413
$ws.Range("J20").Value = 0.007693480963358412
$ws.Range("M20").Value = 8.716382666666666
$ws.Range("N20").Value = 26.149148
$ws.Range("O20").Value = 0.04202276799614693
$ws.Range("P20").Value = 0.04202276799614692
$ws.Range("Q20").Value = 68.02588701409688
$ws.Range("R20").Value = 612.2329831268719
$ws.Range("S20").Value = 0.0003233013656059836
$ws.Range("T20").Value = 0.0003233013656059835

# Row 21
$ws.Range("G21").Value = 7.804371333333333
$ws.Range("H21").Value = 23.413114
$ws.Range("I21").Value = 0.007693480963358413
$ws.Range("J21").Value = 0.007693480963358412
$ws.Range("M21").Value = 36.13945833333333
$ws.Range("N21").Value = 108.418375
$ws.Range("O21").Value = 0.1742328361575779
$ws.Range("P21").Value = 0.1742328361575779
$ws.Range("Q21").Value = 282.0457526188611
$ws.Range("R21").Value = 2538.41177356975
$ws.Range("S21").Value = 0.001340457008170271
$ws.Range("T21").Value = 0.001340457008170271

# Row 22
$ws.Range("G22").Value = 478.3423056666667
$ws.Range("H22").Value = 1435.026917
$ws.Range("I22").Value = 0.4715456588921241
$ws.Range("J22").Value = 0.471545658892124
$ws.Range("M22").Value = 51.402972
$ws.Range("N22").Value = 154.208916
$ws.Range("O22").Value = 0.2478201393026385
$ws.Range("P22").Value = 0.2478201393026385
$ws.Range("Q22").Value = 24588.21614459911
$ws.Range("R22").Value = 221293.945301392
$ws.Range("S22").Value = 0.1168585108742006
$ws.Range("T22").Value = 0.1168585108742006

# Row 23
$ws.Range("G23").Value = 478.3423056666667
$ws.Range("H23").Value = 1435.026917
$ws.Range("I23").Value = 0.4715456588921241
$ws.Range("J23").Value = 0.471545658892124
$ws.Range("O23").Value = 0.2505855219821376
$ws.Range("P23").Value = 0.2505855219821376
$ws.Range("Q23").Value = 24862.5918561026
$ws.Range("R23").Value = 223763.3267049234
$ws.Range("S23").Value = 0.1181625150718939
$ws.Range("T23").Value = 0.1181625150718939

# Row 24
$ws.Range("G24").Value = 478.3423056666667
$ws.Range("H24").Value = 1435.026917
$ws.Range("I24").Value = 0.4715456588921241
$ws.Range("J24").Value = 0.471545658892124
$ws.Range("M24").Value = 59.18509700000001
$ws.Range("N24").Value = 177.555291
$ws.Range("O24").Value = 0.2853387345614992
$ws.Range("P24").Value = 0.2853387345614991
$ws.Range("Q24").Value = 28310.73576008532
$ws.Range("R24").Value = 254796.6218407679
$ws.Range("S24").Value = 0.134550241596247
$ws.Range("T24").Value = 0.134550241596247

# Row 25
$ws.Range("G25").Value = 478.3423056666667
$ws.Range("H25").Value = 1435.026917
$ws.Range("I25").Value = 0.4715456588921241
$ws.Range("J25").Value = 0.471545658892124
$ws.Range("M25").Value = 8.716382666666666
$ws.Range("N25").Value = 26.149148
$ws.Range("O25").Value = 0.04202276799614693
$ws.Range("P25").Value = 0.04202276799614692
$ws.Range("Q25").Value = 4169.414581846302
$ws.Range("R25").Value = 37524.73123661672
$ws.Range("S25").Value = 0.01981565382321397
$ws.Range("T25").Value = 0.01981565382321396

# Row 26
$ws.Range("G26").Value = 478.3423056666667
$ws.Range("H26").Value = 1435.026917
$ws.Range("I26").Value = 0.4715456588921241
$ws.Range("J26").Value = 0.471545658892124
$ws.Range("M26").Value = 36.13945833333333
$ws.Range("N26").Value = 108.418375
$ws.Range("O26").Value = 0.1742328361575779
$ws.Range("P26").Value = 0.1742328361575779
$ws.Range("Q26").Value = 17287.0318247111
$ws.Range("R26").Value = 155583.2864223999
$ws.Range("S26").Value = 0.08215873752656858
$ws.Range("T26").Value = 0.08215873752656856
